$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Value")

$ws.Range("B34").Value = 19134
$ws.Range("B36").Value = 20832
$ws.Range("B37").Value = 31427
$ws.Range("B38").Value = 34296
$ws.Range("B39").Value = 37195
$ws.Range("B40").Value = 37385
$ws.Range("B41").Value = 45842
$ws.Range("B42").Value = 32659
$ws.Range("B43").Value = 38437
$ws.Range("B44").Value = 41526
$ws.Range("B45").Value = 30694
$ws.Range("B46").Value = 47797
$ws.Range("B47").Value = 39648
$ws.Range("B48").Value = 31368
$ws.Range("B49").Value = 45930
$ws.Range("B50").Value = 34668
$ws.Range("B51").Value = 58136
$ws.Range("B52").Value = 35122
$ws.Range("B53").Value = 47139
$ws.Range("B54").Value = 41508
$ws.Range("B55").Value = 59659
$ws.Range("B56").Value = 60009
$ws.Range("B57").Value = 48279
$ws.Range("B58").Value = 48542
$ws.Range("B59").Value = 60998
$ws.Range("B60").Value = 55178
$ws.Range("B61").Value = 73934
$ws.Range("B62").Value = 49558
$ws.Range("B63").Value = 62273
$ws.Range("B64").Value = 56331
$ws.Range("B65").Value = 81767
$ws.Range("B66").Value = 56877
$ws.Range("B67").Value = 69837
$ws.Range("B68").Value = 82904
$ws.Range("B69").Value = 89670
$ws.Range("B70").Value = 64320
$ws.Range("B71").Value = 103334
$ws.Range("B72").Value = 64873
$ws.Range("B73").Value = 65155
$ws.Range("B74").Value = 91603
$ws.Range("B75").Value = 98548
$ws.Range("B76").Value = 72557
$ws.Range("B77").Value = 92704
$ws.Range("B78").Value = 106348
$ws.Range("B79").Value = 113410
$ws.Range("B80").Value = 107122
$ws.Range("B81").Value = 161246
$ws.Range("B82").Value = 94417
$ws.Range("B83").Value = 94765
$ws.Range("B84").Value = 135866
$ws.Range("B85").Value = 149977
$ws.Range("B86").Value = 157331
$ws.Range("B87").Value = 123540
$ws.Range("B88").Value = 137714
$ws.Range("B89").Value = 165783
$ws.Range("B90").Value = 165836
$ws.Range("B91").Value = 188531
$ws.Range("B92").Value = 83453
$ws.Range("B93").Value = 111634
$ws.Range("B94").Value = 193182
$ws.Range("B95").Value = 196591
$ws.Range("B96").Value = 169020
$ws.Range("B97").Value = 169523
$ws.Range("B98").Value = 175684
$ws.Range("B99").Value = 198918
$ws.Range("B100").Value = 170977
$ws.Range("B101").Value = 194304
